{"js": "// The template field placeholder \"<<claimNumber>>\" is renamed to\n// \"<<caseNumber>>\" (the merge-field name changes from \"claimNumber\" to\n// \"caseNumber\"; everything else in the cell, e.g. the \"Case number: \"\n// label and surrounding formatting, stays the same).\nconst body = context.document.body;\nconst results = body.search(\"claimNumber\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"caseNumber\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The template field placeholder \"<<claimNumber>>\" is renamed to\n# \"<<caseNumber>>\" (the merge-field name changes from \"claimNumber\" to\n# \"caseNumber\"; everything else in the cell, e.g. the \"Case number: \"\n# label and surrounding formatting, stays the same).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"claimNumber\"\n$find.Replacement.Text = \"caseNumber\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
